$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous single-row layout (A1:R1)
$ws.Range("A1:R1").Clear()

# New state name / abbreviation pairs (column A = full name, column B = abbreviation)
$states = @(
    @("Colorado", "CO"),
    @("Connecticut", "CT"),
    @("Delaware", "DE"),
    @("Georgia", "GA"),
    @("Idaho", "ID"),
    @("Indiana", "IN"),
    @("Louisiana", "LA"),
    @("Massachusetts", "MA"),
    @("Missouri", "MO"),
    @("New Jersey", "NJ"),
    @("New York", "NY"),
    @("Ohio", "OH"),
    @("Pennsylvania", "PA"),
    @("Rhode Island", "RI"),
    @("Vermont", "VT"),
    @("Washington", "WA")
)

$row = 1
foreach ($pair in $states) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("A1:B16").Select()
